$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New row 24: resale numbers snapshot for 2025-01-14 09:08:59.
# Columns A (Date) and D (Week) look like a date / number to Excel's
# smart-parser, so force them to Text before assigning, then drop the
# number-format override again so the cell keeps the default style.
$ws.Range("A24").NumberFormat = "@"
$ws.Range("D24").NumberFormat = "@"

$ws.Range("A24").Value = "2025-01-14"
$ws.Range("B24").Value = "09:08:59"
$ws.Range("C24").Value = "Tuesday"
$ws.Range("D24").Value = "02"
$ws.Range("E24").Value = 126816
$ws.Range("F24").Value = 143435
$ws.Range("G24").Value = 169128
$ws.Range("H24").Value = 152228
$ws.Range("I24").Value = -1
$ws.Range("J24").Value = 142625
$ws.Range("K24").Value = -1
$ws.Range("L24").Value = -1
$ws.Range("M24").Value = 193041
$ws.Range("N24").Value = 115435
$ws.Range("O24").Value = 45841
$ws.Range("P24").Value = 28485
$ws.Range("Q24").Value = 65312
$ws.Range("R24").Value = -1
$ws.Range("S24").Value = 48136
$ws.Range("T24").Value = -1

$ws.Range("A24").ClearFormats()
$ws.Range("D24").ClearFormats()
